# Update crypto price/volume data per upstream scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.940.29"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.615.36"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'211.36"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.487"
$ws.Range("E7").Value = "  -3.50%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.247"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.0619"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "'18.14"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.840.00"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "1.603.65"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "'4.08"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "25.952.20"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "'61.48"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'191.06"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "'142.92"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("D29").Value = "'15.09"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "'1.21"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("D34").Value = "'2.40"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").Value = "1.121.47"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  -6.38%  "
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("D39").Value = "'0.510"
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").Value = "'97.02"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("D42").Value = "1.751.87"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "'0.750"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").Value = "'5.07"
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'53.71"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -2.10%  "
